$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.36873579987056
$ws.Range("C2").Value = 17.85741523675451
$ws.Range("D2").Value = 6.015478765848312
$ws.Range("E2").Value = 8.693958111024971
$ws.Range("G2").Value = 3.68481227798076
$ws.Range("I2").Value = 44.82148306539234
$ws.Range("M2").Value = 17.48606014150778

$ws.Range("B3").Value = 15.09450607348673
$ws.Range("C3").Value = 17.09403317808741
$ws.Range("D3").Value = 5.904975522039424
$ws.Range("E3").Value = 8.553902495522669
$ws.Range("G3").Value = 3.690957023558072
$ws.Range("I3").Value = 43.23152039762113
$ws.Range("M3").Value = 17.27484477194435

$ws.Range("B4").Value = 14.934965288537
$ws.Range("C4").Value = 16.61446590361236
$ws.Range("D4").Value = 5.838308458828961
$ws.Range("E4").Value = 8.46849814222022
$ws.Range("G4").Value = 3.694906691470082
$ws.Range("I4").Value = 42.22710807657933
$ws.Range("M4").Value = 17.1513682351082

$ws.Range("B5").Value = 14.87228442143178
$ws.Range("C5").Value = 16.4166335358706
$ws.Range("D5").Value = 5.811474768821076
$ws.Range("E5").Value = 8.43387791109563
$ws.Range("G5").Value = 3.696560929062762
$ws.Range("I5").Value = 41.81116714362141
$ws.Range("M5").Value = 17.1026668239254

$ws.Range("B6").Value = 14.86202033078187
$ws.Range("C6").Value = 16.38364854717809
$ws.Range("D6").Value = 5.807040261416971
$ws.Range("E6").Value = 8.428141281401423
$ws.Range("G6").Value = 3.696838322438573
$ws.Range("I6").Value = 41.74171373735396
$ws.Range("M6").Value = 17.09467912278257

$ws.Range("B7").Value = 14.93411036383044
$ws.Range("C7").Value = 16.61180714992028
$ws.Range("D7").Value = 5.837945171542378
$ws.Range("E7").Value = 8.46803045667829
$ws.Range("G7").Value = 3.694928819691341
$ws.Range("I7").Value = 42.22152478989923
$ws.Range("M7").Value = 17.15070481755432

$ws.Range("B8").Value = 15.2724160554235
$ws.Range("C8").Value = 17.59663477724284
$ws.Range("D8").Value = 5.977153774697573
$ws.Range("E8").Value = 8.645565028197725
$ws.Range("G8").Value = 3.686894458708961
$ws.Range("I8").Value = 44.2793655339321
$ws.Range("M8").Value = 17.41197808595826

$ws.Range("B9").Value = 16.00074643102597
$ws.Range("C9").Value = 19.42948256967442
$ws.Range("D9").Value = 6.25791546870614
$ws.Range("E9").Value = 8.996886682853281
$ws.Range("G9").Value = 3.672529482954881
$ws.Range("I9").Value = 48.07385549829333
$ws.Range("M9").Value = 17.97102394163858

$ws.Range("B10").Value = 16.56816552728751
$ws.Range("C10").Value = 20.70203173986706
$ws.Range("D10").Value = 6.466795894093091
$ws.Range("E10").Value = 9.254941976634393
$ws.Range("G10").Value = 3.662805939517407
$ws.Range("I10").Value = 50.69451962184402
$ws.Range("M10").Value = 18.40653911974002

$ws.Range("B11").Value = 16.83174904277735
$ws.Range("C11").Value = 21.2625349374749
$ws.Range("D11").Value = 6.561948673668137
$ws.Range("E11").Value = 9.371913307536403
$ws.Range("G11").Value = 3.658559014999629
$ws.Range("I11").Value = 51.8469572393657
$ws.Range("M11").Value = 18.60922332010001

$ws.Range("B12").Value = 16.93221818971639
$ws.Range("C12").Value = 21.47198164719151
$ws.Range("D12").Value = 6.597964607213042
$ws.Range("E12").Value = 9.41611382525698
$ws.Range("G12").Value = 3.656975878002414
$ws.Range("I12").Value = 52.27740234513711
$ws.Range("M12").Value = 18.6865607203357

$ws.Range("B13").Value = 16.91055321878915
$ws.Range("C13").Value = 21.42700053095709
$ws.Range("D13").Value = 6.590209205306409
$ws.Range("E13").Value = 9.406599156157556
$ws.Range("G13").Value = 3.65731572387354
$ws.Range("I13").Value = 52.1849668540859
$ws.Range("M13").Value = 18.66987983147344

$ws.Range("B14").Value = 16.84000228627333
$ws.Range("C14").Value = 21.27982326269158
$ws.Range("D14").Value = 6.564912250964712
$ws.Range("E14").Value = 9.375551782876384
$ws.Range("G14").Value = 3.658428268316901
$ws.Range("I14").Value = 51.88249092106114
$ws.Range("M14").Value = 18.61557459285064

$ws.Range("B15").Value = 16.79686942140286
$ws.Range("C15").Value = 21.18930349027494
$ws.Range("D15").Value = 6.549413983017038
$ws.Range("E15").Value = 9.356521151803889
$ws.Range("G15").Value = 3.659112991860147
$ws.Range("I15").Value = 51.69643300935414
$ws.Range("M15").Value = 18.58238517178791

$ws.Range("B16").Value = 16.55103847203434
$ws.Range("C16").Value = 20.66501733591094
$ws.Range("D16").Value = 6.460577036956231
$ws.Range("E16").Value = 9.247286279717772
$ws.Range("G16").Value = 3.663087011068136
$ws.Range("I16").Value = 50.61838405471136
$ws.Range("M16").Value = 18.39337901948278

$ws.Range("B17").Value = 16.40153478423084
$ws.Range("C17").Value = 20.33855398504825
$ws.Range("D17").Value = 6.40608664267337
$ws.Range("E17").Value = 9.180142477345994
$ws.Range("G17").Value = 3.665569916053445
$ws.Range("I17").Value = 49.94667732877326
$ws.Range("M17").Value = 18.27854868704598

$ws.Range("B18").Value = 16.31606589229387
$ws.Range("C18").Value = 20.14905821114426
$ws.Range("D18").Value = 6.374759662237215
$ws.Range("E18").Value = 9.141485614640313
$ws.Range("G18").Value = 3.667014634473152
$ws.Range("I18").Value = 49.55660368488984
$ws.Range("M18").Value = 18.21293540496305

$ws.Range("B19").Value = 16.28722119546962
$ws.Range("C19").Value = 20.08460760603208
$ws.Range("D19").Value = 6.3641564688959
$ws.Range("E19").Value = 9.128391666814952
$ws.Range("G19").Value = 3.667506653786985
$ws.Range("I19").Value = 49.42389963741181
$ws.Range("M19").Value = 18.19079660061248

$ws.Range("B20").Value = 16.41739666553986
$ws.Range("C20").Value = 20.37348609439776
$ws.Range("D20").Value = 6.411886003195449
$ws.Range("E20").Value = 9.187294181190479
$ws.Range("G20").Value = 3.665303888786168
$ws.Range("I20").Value = 50.01856914758589
$ws.Range("M20").Value = 18.29072816189166

$ws.Range("B21").Value = 16.86070802202713
$ws.Range("C21").Value = 21.32313003753664
$ws.Range("D21").Value = 6.57234329762221
$ws.Range("E21").Value = 9.384673962093689
$ws.Range("G21").Value = 3.658100808594093
$ws.Range("I21").Value = 51.9714989367962
$ws.Range("M21").Value = 18.63151004139893

$ws.Range("B22").Value = 17.15420282248378
$ws.Range("C22").Value = 21.9273778406689
$ws.Range("D22").Value = 6.67710080681177
$ws.Range("E22").Value = 9.513110521979289
$ws.Range("G22").Value = 3.6535392393747
$ws.Range("I22").Value = 53.21303193288345
$ws.Range("M22").Value = 18.85760930654058

$ws.Range("B23").Value = 16.99725641401127
$ws.Range("C23").Value = 21.6064265491014
$ws.Range("D23").Value = 6.62121113919443
$ws.Range("E23").Value = 9.444623784637516
$ws.Range("G23").Value = 3.655960563439096
$ws.Range("I23").Value = 52.55366152159739
$ws.Range("M23").Value = 18.73665025484066

$ws.Range("B24").Value = 16.41022399540356
$ws.Range("C24").Value = 20.35769890578414
$ws.Range("D24").Value = 6.409264108753583
$ws.Range("E24").Value = 9.184061064613799
$ws.Range("G24").Value = 3.665424105921336
$ws.Range("I24").Value = 49.98607898639257
$ws.Range("M24").Value = 18.28522055689355

$ws.Range("B25").Value = 15.79750515555662
$ws.Range("C25").Value = 18.94576142053448
$ws.Range("D25").Value = 6.18135233847471
$ws.Range("E25").Value = 8.901712882195143
$ws.Range("G25").Value = 3.676268509847844
$ws.Range("I25").Value = 47.07539126635257
$ws.Range("M25").Value = 17.81515917286501
